# Generate Report for Handback
# Adds a new handback row (0502b11e-e274-4c60-a189-f37e77f597c2.md) to the
# Overview sheet and to each per-language detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "0502b11e-e274-4c60-a189-f37e77f597c2.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-11-14 07:13:07"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97629e0edcd0fcdf8af4953ca8f0a22afb2cdf47/e2e/0502b11e-e274-4c60-a189-f37e77f597c2.md",
    "",
    "",
    "e2e\0502b11e-e274-4c60-a189-f37e77f597c2.md"
) | Out-Null

$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "0502b11e-e274-4c60-a189-f37e77f597c2.d46a60c2b093908b1dd5649f21f430558f1bc416.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-11-14 07:12:54"
$wsZhCn.Range("J4").Value = "0502b11e-e274-4c60-a189-f37e77f597c2.d46a60c2b093908b1dd5649f21f430558f1bc416.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-11-14 07:13:38"
$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97629e0edcd0fcdf8af4953ca8f0a22afb2cdf47/e2e/0502b11e-e274-4c60-a189-f37e77f597c2.md",
    "",
    "",
    "0502b11e-e274-4c60-a189-f37e77f597c2.md"
) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d46a60c2b093908b1dd5649f21f430558f1bc416/e2e/0502b11e-e274-4c60-a189-f37e77f597c2.md",
    "",
    "",
    "0502b11e-e274-4c60-a189-f37e77f597c2.md"
) | Out-Null

$tblZhCn = $wsZhCn.ListObjects.Item(1)
$tblZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "0502b11e-e274-4c60-a189-f37e77f597c2.d46a60c2b093908b1dd5649f21f430558f1bc416.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-11-14 07:13:07"
$wsDeDe.Range("J4").Value = "0502b11e-e274-4c60-a189-f37e77f597c2.d46a60c2b093908b1dd5649f21f430558f1bc416.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-11-14 07:13:56"
$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97629e0edcd0fcdf8af4953ca8f0a22afb2cdf47/e2e/0502b11e-e274-4c60-a189-f37e77f597c2.md",
    "",
    "",
    "0502b11e-e274-4c60-a189-f37e77f597c2.md"
) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d46a60c2b093908b1dd5649f21f430558f1bc416/e2e/0502b11e-e274-4c60-a189-f37e77f597c2.md",
    "",
    "",
    "0502b11e-e274-4c60-a189-f37e77f597c2.md"
) | Out-Null

$tblDeDe = $wsDeDe.ListObjects.Item(1)
$tblDeDe.Resize($wsDeDe.Range("A1:P4"))

Write-Output "Report row added to Overview, zh-cn, de-de"
